$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that must stay text (avoid Excel auto-converting numeric-looking strings)
foreach ($addr in @("D5","D6","D10","D13","D14","D17","D20","D22","D23","D24","D25","D28","D31","D32","D35","D38","D39","D40","D41","D47","D49","D51")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "59.306.90"
$ws.Range("E2").Value = "  +1.85%  "
$ws.Range("D3").Value = "2.996.20"
$ws.Range("E3").Value = "  +1.11%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "562.48"
$ws.Range("E5").Value = "  +1.04%  "
$ws.Range("D6").Value = "138.02"
$ws.Range("E6").Value = "  +5.15%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("E8").Value = "  +0.61%  "
$ws.Range("D9").Value = "2.982.98"
$ws.Range("E9").Value = "  +0.72%  "
$ws.Range("D10").Value = "0.133"
$ws.Range("E10").Value = "  +2.57%  "
$ws.Range("E11").Value = "  +5.71%  "
$ws.Range("E12").Value = "  +1.46%  "
$ws.Range("D13").Value = "0.0000230"
$ws.Range("E13").Value = "  +2.55%  "
$ws.Range("D14").Value = "33.66"
$ws.Range("E14").Value = "  +2.27%  "
$ws.Range("E15").Value = "  +1.71%  "
$ws.Range("D16").Value = "3.488.18"
$ws.Range("E16").Value = "  +0.95%  "
$ws.Range("D17").Value = "7.24"
$ws.Range("E17").Value = "  +6.56%  "
$ws.Range("D18").Value = "2.990.45"
$ws.Range("E18").Value = "  +0.71%  "
$ws.Range("D19").Value = "59.227.90"
$ws.Range("E19").Value = "  +1.52%  "
$ws.Range("D20").Value = "429.17"
$ws.Range("E20").Value = "  +2.09%  "
$ws.Range("E21").Value = "  +3.83%  "
$ws.Range("D22").Value = "0.718"
$ws.Range("E22").Value = "  +4.79%  "
$ws.Range("D23").Value = "7.12"
$ws.Range("E23").Value = "  +1.81%  "
$ws.Range("D24").Value = "13.32"
$ws.Range("E24").Value = "  +2.20%  "
$ws.Range("D25").Value = "80.98"
$ws.Range("E25").Value = "  +1.67%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").Value = "2.16"
$ws.Range("E28").Value = "  +7.38%  "
$ws.Range("E29").Value = "  +1.61%  "
$ws.Range("E30").Value = "  +2.45%  "
$ws.Range("D31").Value = "25.73"
$ws.Range("E31").Value = "  +2.30%  "
$ws.Range("D32").Value = "6.12"
$ws.Range("E32").Value = "  -1.22%  "
$ws.Range("E33").Value = "  -5.20%  "
$ws.Range("E34").Value = "  +4.46%  "
$ws.Range("D35").Value = "0.991"
$ws.Range("E35").Value = "  +4.38%  "
$ws.Range("D36").Value = "0.0₃0766"
$ws.Range("E36").Value = "  +12.16%  "
$ws.Range("E37").Value = "  -2.08%  "
$ws.Range("D38").Value = "48.99"
$ws.Range("E38").Value = "  +0.91%  "
$ws.Range("D39").Value = "8.67"
$ws.Range("E39").Value = "  +3.01%  "
$ws.Range("D40").Value = "2.71"
$ws.Range("E40").Value = "  +4.49%  "
$ws.Range("D41").Value = "401.11"
$ws.Range("E41").Value = "  +5.25%  "
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("D43").Value = "2.754.72"
$ws.Range("E43").Value = "  +3.61%  "
$ws.Range("E44").Value = "  -0.99%  "
$ws.Range("E45").Value = "  +4.60%  "
$ws.Range("D47").Value = "34.83"
$ws.Range("E47").Value = "  +20.59%  "
$ws.Range("E48").Value = "  +0.91%  "
$ws.Range("D49").Value = "121.29"
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("E50").Value = "  +0.64%  "
$ws.Range("D51").Value = "23.40"
$ws.Range("E51").Value = "  -0.28%  "
